$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.045.50'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.825.59'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '241.33'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = '0.6361'
$ws.Range("E6").Value = '  -4.37%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '44.87'
$ws.Range("E8").Value = '  +6.89%  '
$ws.Range("D9").Value = '0.2934'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '0.07337'
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '22.78'
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").Value = '0.07668'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.827.55'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Value = '4.986'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '0.6629'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").Value = '81.97'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008693'
$ws.Range("E17").Value = '  +5.05%  '
$ws.Range("D18").Value = '6.036'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").Value = '29.056.46'
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").Value = '2.074.96'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '225.24'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '7.124'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '158.67'
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("D27").Value = '8.469'
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").Value = '0.1366'
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").Value = '17.89'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").Value = '1.503'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").Value = '4.089'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").Value = '4.027'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = '0.05305'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").Value = '  -1.81%  '
$ws.Range("D36").Value = '0.7372'
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("D37").Value = '1.155'
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("D38").Value = '2.651'
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").Value = '1.296.01'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").Value = '0.01785'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").Value = '2.745'
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("D42").Value = '6.298'
$ws.Range("E42").Value = '  +5.74%  '
$ws.Range("D43").Value = '0.8981'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("D45").Value = '102.56'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = '1.973.82'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").Value = '64.01'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").Value = '1.727'
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("D51").Value = '0.07249'
$ws.Range("E51").Value = '  -18.04%  '
